$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price (D) and 1h volume-change (E) columns with latest scrape values.
# Leading apostrophe forces text entry (prevents Excel's numeric auto-detection from
# turning dotted price strings like "1.004" into numbers), and resetting the style
# back to Normal keeps the cell free of an explicit Text number-format.

$ws.Range("D2").Value = "'27.271.37"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.62%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.852.73"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +2.08%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.49%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'314.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +2.04%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'1.003"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.54%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.4644"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.42%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.3719"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +2.24%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.07378"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +2.28%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.8871"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +3.59%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07933"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +5.53%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'20.04"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +1.79%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.834.67"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +1.54%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'5.408"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +1.64%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'6.616"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +1.58%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'92.29"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.68%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'1.005"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.40%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.000008948"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +4.43%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  -0.38%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'14.90"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +3.44%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'27.290.79"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +1.81%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'5.145"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.06%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'10.59"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.86%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'2.075.03"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +2.46%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'153.06"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +1.20%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'1.865"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +1.36%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'18.54"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +2.60%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +0.37%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'5.148"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +1.36%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'117.30"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +2.14%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.08896"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +0.46%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'0.7499"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +4.77%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'2.972"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.68%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'4.488"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +1.97%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.147"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +1.58%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'2.576"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +6.32%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'1.082"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.89%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.05288"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +1.05%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.01953"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +2.11%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'2.980"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +2.07%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'7.141"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +0.07%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.5188"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +1.25%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.1640"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +1.25%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'8.318"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +1.86%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.4883"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +2.00%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'10.25"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +1.88%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'1.003"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.51%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'103.30"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  +0.58%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'1.638"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +1.60%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.06245"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +0.96%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'65.44"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +1.97%  "
$ws.Range("E51").Style = "Normal"
